$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("H2").Value = 0.725
$ws.Range("E3").Value = 0.9999420137639272
$ws.Range("G3").Value = 0.7321428571428572
$ws.Range("I3").Value = 0.9939823107496032
$ws.Range("F4").Value = 0.98
$ws.Range("G4").Value = 0.6285714285714286
$ws.Range("I4").Value = 0.8028441704858906
$ws.Range("F5").Value = 0.68
$ws.Range("G5").Value = 0.6799999999999999
$ws.Range("I5").Value = 0.9758384973814157
$ws.Range("E6").Value = 0.95
$ws.Range("G6").Value = 0.54
$ws.Range("I6").Value = 0.6567941620385526
$ws.Range("E7").Value = 0.9999464535798339
$ws.Range("F7").Value = 0.97
$ws.Range("I7").Value = 0.9735480799325416
$ws.Range("E8").Value = 0.7580255003625727
$ws.Range("F8").Value = 0.9399999999999999
$ws.Range("G8").Value = 0.7414285714285714
$ws.Range("I8").Value = 0.9839499946097089
$ws.Range("E9").Value = 0.6034027588324409
$ws.Range("F9").Value = 0.71
$ws.Range("I9").Value = 0.8871990296474435
$ws.Range("E11").Value = 0.9327235094660271
$ws.Range("F11").Value = 1
$ws.Range("G11").Value = 0.3078571428571428
$ws.Range("I11").Value = 0.3205809741773398
$ws.Range("E12").Value = 0.999999958096461
$ws.Range("F12").Value = 0.98
$ws.Range("G12").Value = 0.3307142857142857
$ws.Range("I12").Value = 0.4325513506965651
$ws.Range("F13").Value = 0.7
$ws.Range("G13").Value = 0.7
$ws.Range("I13").Value = 0.9585537977361916
$ws.Range("E14").Value = 0.8643225808107002
$ws.Range("G14").Value = 0.6978571428571428
$ws.Range("I14").Value = 0.8902288467931158
$ws.Range("E15").Value = 0.8163265148939988
$ws.Range("I15").Value = 0.8932822971611728
$ws.Range("E16").Value = 0.9996142192362972
$ws.Range("F16").Value = 0.9399999999999999
$ws.Range("G16").Value = 0.6342857142857143
$ws.Range("I16").Value = 0.8553813521227718
$ws.Range("E17").Value = 0.9999985611864662
$ws.Range("F17").Value = 0.98
$ws.Range("G17").Value = 0.6528571428571428
$ws.Range("I17").Value = 0.9838503259939992
$ws.Range("E18").Value = 0.9915681842605448
$ws.Range("G18").Value = 0.8571428571428572
$ws.Range("I18").Value = 0.8567978043241556
$ws.Range("E19").Value = 0.9999981868198751
$ws.Range("F19").Value = 0.96
$ws.Range("G19").Value = 0.7007142857142857
$ws.Range("I19").Value = 0.9174374620405917
$ws.Range("E20").Value = 0.9598807654112183
$ws.Range("F20").Value = 1
$ws.Range("G20").Value = 0.895
$ws.Range("I20").Value = 0.9509061379344278
$ws.Range("E21").Value = 0.8374922442792085
$ws.Range("F21").Value = 0.99
$ws.Range("G21").Value = 0.6942857142857143
$ws.Range("I21").Value = 0.9723038473978806
